$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.035907263186886
$ws.Range("D2").Value = 1.038331459212118
$ws.Range("E2").Value = 1.054705496752287
$ws.Range("F2").Value = 1.060881494940947
$ws.Range("I2").Value = 1.038843777116449
$ws.Range("J2").Value = 1.04101847607822
$ws.Range("K2").Value = 1.04111979618313
$ws.Range("L2").Value = 1.057447970662354
$ws.Range("M2").Value = 1.063607068352655
$ws.Range("N2").Value = 1.017605315334003
# Row 3
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.036816276142956
$ws.Range("D3").Value = 1.038997993348583
$ws.Range("E3").Value = 1.055792441743521
$ws.Range("F3").Value = 1.061993703966865
$ws.Range("I3").Value = 1.039069726821165
$ws.Range("J3").Value = 1.041571436739696
$ws.Range("K3").Value = 1.041596700931714
$ws.Range("L3").Value = 1.058347512855038
$ws.Range("M3").Value = 1.064533038625436
$ws.Range("N3").Value = 1.017790803211627
# Row 4
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.037404681548745
$ws.Range("D4").Value = 1.039429462923815
$ws.Range("E4").Value = 1.056496857419016
$ws.Range("F4").Value = 1.062714345452873
$ws.Range("I4").Value = 1.03921489065121
$ws.Range("J4").Value = 1.041928820224057
$ws.Range("K4").Value = 1.041904785686857
$ws.Range("L4").Value = 1.058930069327272
$ws.Range("M4").Value = 1.065132577976236
$ws.Range("N4").Value = 1.01791063160903
# Row 5
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.037652097373906
$ws.Range("D5").Value = 1.039610893991545
$ws.Range("E5").Value = 1.056793254128261
$ws.Range("F5").Value = 1.063017534314977
$ws.Range("I5").Value = 1.039275667926345
$ws.Range("J5").Value = 1.042078963142856
$ws.Range("K5").Value = 1.042034182940321
$ws.Range("L5").Value = 1.059175093527305
$ws.Range("M5").Value = 1.065384713492624
$ws.Range("N5").Value = 1.017960960541901
# Row 6
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.037693642498107
$ws.Range("D6").Value = 1.039641359450179
$ws.Range("E6").Value = 1.056843035693467
$ws.Range("F6").Value = 1.06306845460036
$ws.Range("I6").Value = 1.039285858048791
$ws.Range("J6").Value = 1.042104166848523
$ws.Range("K6").Value = 1.042055902156011
$ws.Range("L6").Value = 1.059216241043503
$ws.Range("M6").Value = 1.065427053369383
$ws.Range("N6").Value = 1.017969408230076
# Row 7
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.037407987334951
$ws.Range("D7").Value = 1.039431887054032
$ws.Range("E7").Value = 1.056500816863916
$ws.Range("F7").Value = 1.062718395769552
$ws.Range("I7").Value = 1.039215703741439
$ws.Range("J7").Value = 1.041930826839305
$ws.Range("K7").Value = 1.041906515178006
$ws.Range("L7").Value = 1.058933342891645
$ws.Range("M7").Value = 1.065135946673599
$ws.Range("N7").Value = 1.017911304291483
# Row 8
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.036214423702497
$ws.Range("D8").Value = 1.03855668007773
$ws.Range("E8").Value = 1.055072608999507
$ws.Range("F8").Value = 1.061257170463323
$ws.Range("I8").Value = 1.038920353101081
$ws.Range("J8").Value = 1.041205438018786
$ws.Range("K8").Value = 1.041281072357646
$ws.Range("L8").Value = 1.057751872904646
$ws.Range("M8").Value = 1.063919926448757
$ws.Range("N8").Value = 1.017668041942899
# Row 9
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.034112883006433
$ws.Range("D9").Value = 1.037015865977025
$ws.Range("E9").Value = 1.052564290736967
$ws.Range("F9").Value = 1.058689740884968
$ws.Range("I9").Value = 1.038391956347688
$ws.Range("J9").Value = 1.039924035228502
$ws.Range("K9").Value = 1.040175133426568
$ws.Range("L9").Value = 1.055673763881989
$ws.Range("M9").Value = 1.06178003597697
$ws.Range("N9").Value = 1.01723790593172
# Row 10
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.032713028704363
$ws.Range("D10").Value = 1.035989677752668
$ws.Range("E10").Value = 1.050897730856612
$ws.Range("F10").Value = 1.056983150274158
$ws.Range("I10").Value = 1.038034371886162
$ws.Range("J10").Value = 1.039067677376794
$ws.Range("K10").Value = 1.039435313207637
$ws.Range("L10").Value = 1.054290932603561
$ws.Range("M10").Value = 1.060355416455629
$ws.Range("N10").Value = 1.016950174622368
# Row 11
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.032107164432507
$ws.Range("D11").Value = 1.035545583994066
$ws.Range("E11").Value = 1.050177437774873
$ws.Range("F11").Value = 1.056245376102465
$ws.Range("I11").Value = 1.037878277347312
$ws.Range("J11").Value = 1.038696378675515
$ws.Range("K11").Value = 1.039114372963054
$ws.Range("L11").Value = 1.053692766843949
$ws.Range("M11").Value = 1.059739014334703
$ws.Range("N11").Value = 1.016825356655463
# Row 12
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.031882162550099
$ws.Range("D12").Value = 1.03538066690742
$ws.Range("E12").Value = 1.049910090262043
$ws.Range("F12").Value = 1.05597151320477
$ws.Range("I12").Value = 1.037820108261565
$ws.Range("J12").Value = 1.038558388875827
$ws.Range("K12").Value = 1.038995073078883
$ws.Range("L12").Value = 1.053470673303561
$ws.Range("M12").Value = 1.059510125763966
$ws.Range("N12").Value = 1.016778959686866
# Row 13
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.031930424256169
$ws.Range("D13").Value = 1.0354160403938
$ws.Range("E13").Value = 1.049967428091417
$ws.Range("F13").Value = 1.056030249607192
$ws.Range("I13").Value = 1.037832594257224
$ws.Range("J13").Value = 1.038587991445227
$ws.Range("K13").Value = 1.039020667296919
$ws.Range("L13").Value = 1.053518308960101
$ws.Range("M13").Value = 1.059559219937386
$ws.Range("N13").Value = 1.016788913528739
# Row 14
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.032088564823545
$ws.Range("D14").Value = 1.035531951091667
$ws.Range("E14").Value = 1.050155334635407
$ws.Range("F14").Value = 1.056222734850499
$ws.Range("I14").Value = 1.037873472918131
$ws.Range("J14").Value = 1.038684973879309
$ws.Range("K14").Value = 1.039104513396235
$ws.Range("L14").Value = 1.053674406646942
$ws.Range("M14").Value = 1.059720092891145
$ws.Range("N14").Value = 1.016821522159802
# Row 15
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.032186006177546
$ws.Range("D15").Value = 1.035603372715428
$ws.Range("E15").Value = 1.050271136757696
$ws.Range("F15").Value = 1.056341355129157
$ws.Range("I15").Value = 1.037898634627562
$ws.Range("J15").Value = 1.038744718314872
$ws.Range("K15").Value = 1.039156162055103
$ws.Range("L15").Value = 1.053770595765007
$ws.Range("M15").Value = 1.059819221397072
$ws.Range("N15").Value = 1.016841608914796
# Row 16
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.032753243922248
$ws.Range("D16").Value = 1.036019156213033
$ws.Range("E16").Value = 1.050945562562265
$ws.Range("F16").Value = 1.057032139036959
$ws.Range("I16").Value = 1.0380447048932
$ws.Range("J16").Value = 1.039092308982801
$ws.Range("K16").Value = 1.039456600519224
$ws.Range("L16").Value = 1.054330643759144
$ws.Range("M16").Value = 1.060396334917061
$ws.Range("N16").Value = 1.016958453592834
# Row 17
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.033109133158124
$ws.Range("D17").Value = 1.036280034648819
$ws.Range("E17").Value = 1.051368970600548
$ws.Range("F17").Value = 1.057465768609586
$ws.Range("I17").Value = 1.038135994357103
$ws.Range("J17").Value = 1.039310212707047
$ws.Range("K17").Value = 1.039644899330666
$ws.Range("L17").Value = 1.054682110419027
$ws.Range("M17").Value = 1.060758468419793
$ws.Range("N17").Value = 1.017031686230657
# Row 18
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.033316744475955
$ws.Range("D18").Value = 1.036432224912237
$ws.Range("E18").Value = 1.051616066453581
$ws.Range("F18").Value = 1.057718812274915
$ws.Range("I18").Value = 1.038189120541358
$ws.Range("J18").Value = 1.039437264928033
$ws.Range("K18").Value = 1.039754673556339
$ws.Range("L18").Value = 1.054887173743499
$ws.Range("M18").Value = 1.060969739849396
$ws.Range("N18").Value = 1.017074379530956
# Row 19
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.033387539152582
$ws.Range("D19").Value = 1.036484121942704
$ws.Range("E19").Value = 1.051700341642678
$ws.Range("F19").Value = 1.057805113098659
$ws.Range("I19").Value = 1.038207214591192
$ws.Range("J19").Value = 1.039480578401932
$ws.Range("K19").Value = 1.039792093991533
$ws.Range("L19").Value = 1.054957104996451
$ws.Range("M19").Value = 1.061041785582729
$ws.Range("N19").Value = 1.017088933092515
# Row 20
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.033070946797841
$ws.Range("D20").Value = 1.036252042341139
$ws.Range("E20").Value = 1.051323529568401
$ws.Range("F20").Value = 1.057419232375766
$ws.Range("I20").Value = 1.038126212417312
$ws.Range("J20").Value = 1.039286838590607
$ws.Range("K20").Value = 1.03962470257019
$ws.Range("L20").Value = 1.054644395291678
$ws.Range("M20").Value = 1.060719610256058
$ws.Range("N20").Value = 1.017023831341122
# Row 21
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.032041995172743
$ws.Range("D21").Value = 1.035497817191158
$ws.Range("E21").Value = 1.050099995279036
$ws.Range("F21").Value = 1.056166047791918
$ws.Range("I21").Value = 1.037861440375148
$ws.Range("J21").Value = 1.038656416977274
$ws.Range("K21").Value = 1.039079825243814
$ws.Range("L21").Value = 1.053628437230569
$ws.Range("M21").Value = 1.059672717863306
$ws.Range("N21").Value = 1.016811920667377
# Row 22
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.031395302086075
$ws.Range("D22").Value = 1.035023832752957
$ws.Range("E22").Value = 1.04933187617053
$ws.Range("F22").Value = 1.055379158347524
$ws.Range("I22").Value = 1.03769387655955
$ws.Range("J22").Value = 1.038259624285951
$ws.Range("K22").Value = 1.038736728224435
$ws.Range("L22").Value = 1.052990195685192
$ws.Range("M22").Value = 1.059014903856776
$ws.Range("N22").Value = 1.016678487344919
# Row 23
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.0317381024178
$ws.Range("D23").Value = 1.035275078978475
$ws.Range("E23").Value = 1.049738960056504
$ws.Range("F23").Value = 1.055796204867184
$ws.Range("I23").Value = 1.037782808602942
$ws.Range("J23").Value = 1.038470011327365
$ws.Range("K23").Value = 1.038918658638959
$ws.Range("L23").Value = 1.053328489107249
$ws.Range("M23").Value = 1.05936358463039
$ws.Range("N23").Value = 1.016749241416465
# Row 24
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.033088201486478
$ws.Range("D24").Value = 1.036264690785141
$ws.Range("E24").Value = 1.051344062015863
$ws.Range("F24").Value = 1.057440259741707
$ws.Range("I24").Value = 1.038130632830021
$ws.Range("J24").Value = 1.039297400493997
$ws.Range("K24").Value = 1.039633828793946
$ws.Range("L24").Value = 1.054661436952664
$ws.Range("M24").Value = 1.060737168448272
$ws.Range("N24").Value = 1.017027380695773
# Row 25
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.034655978558912
$ws.Range("D25").Value = 1.037414028597893
$ws.Range("E25").Value = 1.053211756886982
$ws.Range("F25").Value = 1.059352598739941
$ws.Range("I25").Value = 1.03852949916174
$ws.Range("J25").Value = 1.040255680148291
$ws.Range("K25").Value = 1.04046149407286
$ws.Range("L25").Value = 1.056210552925497
$ws.Range("M25").Value = 1.062332903173696
$ws.Range("N25").Value = 1.017349279449784
